$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "684÷5="
$t.Cell(1, 2).Range.Text = "423÷6="
$t.Cell(1, 3).Range.Text = "867÷7="
$t.Cell(1, 4).Range.Text = "402÷4="
$t.Cell(1, 5).Range.Text = "707÷5="

$t.Cell(5, 1).Range.Text = "200÷4="
$t.Cell(5, 2).Range.Text = "754÷6="
$t.Cell(5, 3).Range.Text = "653÷7="
$t.Cell(5, 4).Range.Text = "726÷9="
$t.Cell(5, 5).Range.Text = "985÷5="

$t.Cell(9, 1).Range.Text = "721÷3="
$t.Cell(9, 2).Range.Text = "978÷3="
$t.Cell(9, 3).Range.Text = "549÷4="
$t.Cell(9, 4).Range.Text = "218÷9="
$t.Cell(9, 5).Range.Text = "755÷8="

$t.Cell(13, 1).Range.Text = "856÷7="
$t.Cell(13, 2).Range.Text = "372÷3="
$t.Cell(13, 3).Range.Text = "989÷5="
$t.Cell(13, 4).Range.Text = "899÷8="
$t.Cell(13, 5).Range.Text = "915÷7="

$t.Cell(17, 1).Range.Text = "389÷6="
$t.Cell(17, 2).Range.Text = "376÷5="
$t.Cell(17, 3).Range.Text = "860÷3="
$t.Cell(17, 4).Range.Text = "896÷9="
$t.Cell(17, 5).Range.Text = "594÷9="

Write-Output "Updated 25 cells."